# Auto-generated edit script: refreshes market-price-derived columns (H-N)
# on the Raiden_Profits workbook per the scheduled runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (Leve Item ID 4564)
$ws.Cells.Item(6, 8).Value = 245.2
$ws.Cells.Item(6, 9).Value = 245.2
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 735.5999999999999
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -623.5999999999999
$ws.Cells.Item(6, 14).ClearContents()
# Row 9 (Leve Item ID 5487)
$ws.Cells.Item(9, 8).Value = 555814.7
$ws.Cells.Item(9, 9).Value = 258
$ws.Cells.Item(9, 10).Value = 1250260.5
$ws.Cells.Item(9, 11).Value = 258
$ws.Cells.Item(9, 12).Value = 1250260.5
$ws.Cells.Item(9, 13).Value = -89
# Row 11 (Leve Item ID 5533)
$ws.Cells.Item(11, 8).Value = 744.1667
$ws.Cells.Item(11, 9).Value = 744.1667
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 744.1667
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -604.1667
# Row 17 (Leve Item ID 38956)
$ws.Cells.Item(17, 8).Value = 411.64865
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 411.64865
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1234.94595
$ws.Cells.Item(17, 14).Value = -1570.94595
# Row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 2699.5
$ws.Cells.Item(40, 9).Value = 2101
$ws.Cells.Item(40, 10).Value = 2998.75
$ws.Cells.Item(40, 11).Value = 2101
$ws.Cells.Item(40, 12).Value = 2998.75
$ws.Cells.Item(40, 13).Value = -1926
$ws.Cells.Item(40, 14).Value = -3348.75
# Row 51 (Leve Item ID 5486)
$ws.Cells.Item(51, 8).Value = 7245.952
$ws.Cells.Item(51, 9).Value = 5815.5
$ws.Cells.Item(51, 10).Value = 7818.1333
$ws.Cells.Item(51, 11).Value = 5815.5
$ws.Cells.Item(51, 12).Value = 7818.1333
$ws.Cells.Item(51, 13).Value = -5331.5
# Row 53 (Leve Item ID 5479)
$ws.Cells.Item(53, 8).Value = 1198.9286
$ws.Cells.Item(53, 9).Value = 708.1818
$ws.Cells.Item(53, 10).Value = 2998.3333
$ws.Cells.Item(53, 11).Value = 708.1818
$ws.Cells.Item(53, 12).Value = 2998.3333
$ws.Cells.Item(53, 13).Value = -71.18179999999995
$ws.Cells.Item(53, 14).Value = -4272.3333
# Row 64 (Leve Item ID 5506)
$ws.Cells.Item(64, 8).Value = 9997.5
$ws.Cells.Item(64, 9).Value = 9998
$ws.Cells.Item(64, 10).Value = 9997
$ws.Cells.Item(64, 11).Value = 9998
$ws.Cells.Item(64, 12).Value = 9997
$ws.Cells.Item(64, 13).Value = -9750
$ws.Cells.Item(64, 14).Value = -10493
# Row 67 (Leve Item ID 5506)
$ws.Cells.Item(67, 8).Value = 9997.5
$ws.Cells.Item(67, 9).Value = 9998
$ws.Cells.Item(67, 10).Value = 9997
$ws.Cells.Item(67, 11).Value = 9998
$ws.Cells.Item(67, 12).Value = 9997
$ws.Cells.Item(67, 13).Value = -9140
$ws.Cells.Item(67, 14).Value = -11713
# Row 86 (Leve Item ID 12603)
$ws.Cells.Item(86, 8).Value = 1714.84
$ws.Cells.Item(86, 9).Value = 1998.5454
$ws.Cells.Item(86, 10).Value = 1491.9286
$ws.Cells.Item(86, 11).Value = 1998.5454
$ws.Cells.Item(86, 12).Value = 1491.9286
$ws.Cells.Item(86, 13).Value = -875.5454
$ws.Cells.Item(86, 14).Value = -3737.9286
# Row 89 (Leve Item ID 12603)
$ws.Cells.Item(89, 8).Value = 1714.84
$ws.Cells.Item(89, 9).Value = 1998.5454
$ws.Cells.Item(89, 10).Value = 1491.9286
$ws.Cells.Item(89, 11).Value = 9992.726999999999
$ws.Cells.Item(89, 12).Value = 7459.643
$ws.Cells.Item(89, 13).Value = -4376.726999999999
$ws.Cells.Item(89, 14).Value = -18691.643
# Row 106 (Leve Item ID 19903)
$ws.Cells.Item(106, 8).Value = 1649.5
$ws.Cells.Item(106, 9).Value = 1399.6666
$ws.Cells.Item(106, 10).Value = 2399
$ws.Cells.Item(106, 11).Value = 1399.6666
$ws.Cells.Item(106, 12).Value = 2399
$ws.Cells.Item(106, 13).Value = -768.6666
$ws.Cells.Item(106, 14).Value = -3661
# Row 107 (Leve Item ID 27766)
$ws.Cells.Item(107, 8).Value = 1928.9286
$ws.Cells.Item(107, 9).Value = 1357
$ws.Cells.Item(107, 10).Value = 3358.75
$ws.Cells.Item(107, 11).Value = 1357
$ws.Cells.Item(107, 12).Value = 3358.75
$ws.Cells.Item(107, 13).Value = 563
$ws.Cells.Item(107, 14).Value = -7198.75
# Row 126 (Leve Item ID 34391)
$ws.Cells.Item(126, 8).Value = 78317
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 78317
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 78317
$ws.Cells.Item(126, 14).Value = -88197
# Row 132 (Leve Item ID 44049)
$ws.Cells.Item(132, 8).Value = 3409
$ws.Cells.Item(132, 9).Value = 2986.5881
$ws.Cells.Item(132, 10).Value = 6999.5
$ws.Cells.Item(132, 11).Value = 8959.764299999999
$ws.Cells.Item(132, 12).Value = 20998.5
$ws.Cells.Item(132, 13).Value = -6429.764299999999
$ws.Cells.Item(132, 14).Value = -26058.5
# Row 137 (Leve Item ID 44013)
$ws.Cells.Item(137, 8).Value = 1823.6471
$ws.Cells.Item(137, 9).Value = 1213
$ws.Cells.Item(137, 10).Value = 2251.1
$ws.Cells.Item(137, 11).Value = 3639
$ws.Cells.Item(137, 12).Value = 6753.299999999999
$ws.Cells.Item(137, 13).Value = -1089
# Row 138 (Leve Item ID 44169)
$ws.Cells.Item(138, 8).Value = 2496.3696
$ws.Cells.Item(138, 9).Value = 3207.2
$ws.Cells.Item(138, 10).Value = 2298.9167
$ws.Cells.Item(138, 11).Value = 9621.599999999999
$ws.Cells.Item(138, 12).Value = 6896.750100000001
$ws.Cells.Item(138, 13).Value = -4481.599999999999
$ws.Cells.Item(138, 14).Value = -17176.7501
# Row 141 (Leve Item ID 44161)
$ws.Cells.Item(141, 8).Value = 5371
$ws.Cells.Item(141, 9).Value = 3477
$ws.Cells.Item(141, 10).Value = 12000
$ws.Cells.Item(141, 11).Value = 10431
$ws.Cells.Item(141, 12).Value = 36000
$ws.Cells.Item(141, 13).Value = -5251
$ws.Cells.Item(141, 14).Value = -46360

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID 27714)
$ws.Cells.Item(45, 8).Value = 1740.5
$ws.Cells.Item(45, 9).Value = 1569.4
$ws.Cells.Item(45, 10).Value = 1826.05
$ws.Cells.Item(45, 11).Value = 1569.4
$ws.Cells.Item(45, 12).Value = 1826.05
$ws.Cells.Item(45, 13).Value = -1192.4
$ws.Cells.Item(45, 14).Value = -2580.05
# Row 63 (Leve Item ID 12528)
$ws.Cells.Item(63, 8).Value = 4032.4
$ws.Cells.Item(63, 9).Value = 2295.5
$ws.Cells.Item(63, 10).Value = 5190.3335
$ws.Cells.Item(63, 11).Value = 2295.5
$ws.Cells.Item(63, 12).Value = 5190.3335
$ws.Cells.Item(63, 13).Value = -1609.5
$ws.Cells.Item(63, 14).Value = -6562.3335
# Row 66 (Leve Item ID 12528)
$ws.Cells.Item(66, 8).Value = 4032.4
$ws.Cells.Item(66, 9).Value = 2295.5
$ws.Cells.Item(66, 10).Value = 5190.3335
$ws.Cells.Item(66, 11).Value = 11477.5
$ws.Cells.Item(66, 12).Value = 25951.6675
$ws.Cells.Item(66, 13).Value = -8045.5
$ws.Cells.Item(66, 14).Value = -32815.6675
# Row 102 (Leve Item ID 19945)
$ws.Cells.Item(102, 8).Value = 2910
$ws.Cells.Item(102, 9).Value = 2910
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 2910
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -1288
# Row 110 (Leve Item ID 27708)
$ws.Cells.Item(110, 8).Value = 1328.3334
$ws.Cells.Item(110, 9).Value = 1257.7
$ws.Cells.Item(110, 10).Value = 1469.6
$ws.Cells.Item(110, 11).Value = 1257.7
$ws.Cells.Item(110, 12).Value = 1469.6
$ws.Cells.Item(110, 13).Value = 787.3
$ws.Cells.Item(110, 14).Value = -5559.6
# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 3537.4167
$ws.Cells.Item(132, 9).Value = 3208.2856
$ws.Cells.Item(132, 10).Value = 3998.2
$ws.Cells.Item(132, 11).Value = 9624.856800000001
$ws.Cells.Item(132, 12).Value = 11994.6
$ws.Cells.Item(132, 13).Value = -7094.856800000001
$ws.Cells.Item(132, 14).Value = -17054.6

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 2897.8667
$ws.Cells.Item(134, 9).Value = 2006.8
$ws.Cells.Item(134, 10).Value = 4680
$ws.Cells.Item(134, 11).Value = 6020.4
$ws.Cells.Item(134, 12).Value = 14040
$ws.Cells.Item(134, 13).Value = -3485.4
$ws.Cells.Item(134, 14).Value = -19110

$ws = $wb.Worksheets.Item("CRP")
# Row 3 (Leve Item ID 3763)
$ws.Cells.Item(3, 8).Value = 4118.6
$ws.Cells.Item(3, 9).Value = 198
$ws.Cells.Item(3, 10).Value = 9999.5
$ws.Cells.Item(3, 11).Value = 198
$ws.Cells.Item(3, 12).Value = 9999.5
$ws.Cells.Item(3, 13).Value = -85
# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 3025.5833
$ws.Cells.Item(31, 9).Value = 2056.5
$ws.Cells.Item(31, 10).Value = 3510.125
$ws.Cells.Item(31, 11).Value = 2056.5
$ws.Cells.Item(31, 12).Value = 3510.125
$ws.Cells.Item(31, 13).Value = -1761.5
$ws.Cells.Item(31, 14).Value = -4100.125
# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 3025.5833
$ws.Cells.Item(34, 9).Value = 2056.5
$ws.Cells.Item(34, 10).Value = 3510.125
$ws.Cells.Item(34, 11).Value = 2056.5
$ws.Cells.Item(34, 12).Value = 3510.125
$ws.Cells.Item(34, 13).Value = -1854.5
$ws.Cells.Item(34, 14).Value = -3914.125
# Row 103 (Leve Item ID 19558)
$ws.Cells.Item(103, 8).Value = 3500
$ws.Cells.Item(103, 9).Value = 3500
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 3500
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = -2328

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Cells.Item(4, 8).Value = 77155700
$ws.Cells.Item(4, 9).Value = 3879898.5
$ws.Cells.Item(4, 10).Value = 700000000
$ws.Cells.Item(4, 11).Value = 11639695.5
$ws.Cells.Item(4, 12).Value = 2100000000
$ws.Cells.Item(4, 13).Value = -11639583.5
# Row 37 (Leve Item ID 9516)
$ws.Cells.Item(37, 8).Value = 79906.75
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 79906.75
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 239720.25
$ws.Cells.Item(37, 14).Value = -239944.25
# Row 113 (Leve Item ID 27843)
$ws.Cells.Item(113, 8).Value = 549.36365
$ws.Cells.Item(113, 9).Value = 400.14285
$ws.Cells.Item(113, 10).Value = 810.5
$ws.Cells.Item(113, 11).Value = 1200.42855
$ws.Cells.Item(113, 12).Value = 2431.5
$ws.Cells.Item(113, 13).Value = 969.5714499999999
$ws.Cells.Item(113, 14).Value = -6771.5
# Row 120 (Leve Item ID 27877)
$ws.Cells.Item(120, 8).Value = 10399.8
$ws.Cells.Item(120, 9).Value = 3999.6667
$ws.Cells.Item(120, 10).Value = 20000
$ws.Cells.Item(120, 11).Value = 11999.0001
$ws.Cells.Item(120, 12).Value = 60000
$ws.Cells.Item(120, 13).Value = -7161.000100000001
# Row 140 (Leve Item ID 44097)
$ws.Cells.Item(140, 8).Value = 2142.611
$ws.Cells.Item(140, 9).Value = 2142.611
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 6427.833
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -1247.833
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Cells.Item(70, 8).Value = 9224
$ws.Cells.Item(70, 9).Value = 5989
$ws.Cells.Item(70, 10).Value = 10518
$ws.Cells.Item(70, 11).Value = 5989
$ws.Cells.Item(70, 12).Value = 10518
$ws.Cells.Item(70, 13).Value = -5719
$ws.Cells.Item(70, 14).Value = -11058
# Row 73 (Leve Item ID 14146)
$ws.Cells.Item(73, 8).Value = 9224
$ws.Cells.Item(73, 9).Value = 5989
$ws.Cells.Item(73, 10).Value = 10518
$ws.Cells.Item(73, 11).Value = 5989
$ws.Cells.Item(73, 12).Value = 10518
$ws.Cells.Item(73, 13).Value = -5053
$ws.Cells.Item(73, 14).Value = -12390
# Row 95 (Leve Item ID 18235)
$ws.Cells.Item(95, 8).Value = 54361.25
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 54361.25
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 54361.25
$ws.Cells.Item(95, 14).Value = -59853.25
# Row 102 (Leve Item ID 36169)
$ws.Cells.Item(102, 8).Value = 3975.818
$ws.Cells.Item(102, 9).Value = 3565.75
$ws.Cells.Item(102, 10).Value = 4210.143
$ws.Cells.Item(102, 11).Value = 3565.75
$ws.Cells.Item(102, 12).Value = 4210.143
$ws.Cells.Item(102, 13).Value = -1943.75
# Row 132 (Leve Item ID 44008)
$ws.Cells.Item(132, 8).Value = 5240.4443
$ws.Cells.Item(132, 9).Value = 4880.5713
$ws.Cells.Item(132, 10).Value = 6500
$ws.Cells.Item(132, 11).Value = 14641.7139
$ws.Cells.Item(132, 12).Value = 19500
$ws.Cells.Item(132, 13).Value = -12111.7139

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Cells.Item(22, 8).Value = 1350
$ws.Cells.Item(22, 9).Value = 623.6667
$ws.Cells.Item(22, 10).Value = 1785.8
$ws.Cells.Item(22, 11).Value = 623.6667
$ws.Cells.Item(22, 12).Value = 1785.8
$ws.Cells.Item(22, 13).Value = -328.6667
# Row 27 (Leve Item ID 5277)
$ws.Cells.Item(27, 8).Value = 1350
$ws.Cells.Item(27, 9).Value = 623.6667
$ws.Cells.Item(27, 10).Value = 1785.8
$ws.Cells.Item(27, 11).Value = 623.6667
$ws.Cells.Item(27, 12).Value = 1785.8
$ws.Cells.Item(27, 13).Value = -516.6667
# Row 55 (Leve Item ID 5284)
$ws.Cells.Item(55, 8).Value = 555.2143
$ws.Cells.Item(55, 9).Value = 413.6
$ws.Cells.Item(55, 10).Value = 909.25
$ws.Cells.Item(55, 11).Value = 413.6
$ws.Cells.Item(55, 12).Value = 909.25
$ws.Cells.Item(55, 13).Value = -240.6
# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 5085.857
$ws.Cells.Item(132, 9).Value = 2502
$ws.Cells.Item(132, 10).Value = 6119.4
$ws.Cells.Item(132, 11).Value = 7506
$ws.Cells.Item(132, 12).Value = 18358.2
$ws.Cells.Item(132, 13).Value = -4976
$ws.Cells.Item(132, 14).Value = -23418.2

$ws = $wb.Worksheets.Item("WVR")
# Row 7 (Leve Item ID 2661)
$ws.Cells.Item(7, 8).Value = 6671000
$ws.Cells.Item(7, 9).Value = 10004000
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 11).Value = 10004000
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = -10003887
$ws.Cells.Item(7, 14).Value = -5226
# Row 9 (Leve Item ID 3015)
$ws.Cells.Item(9, 8).Value = 480
$ws.Cells.Item(9, 9).Value = 480
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 480
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = -340
# Row 81 (Leve Item ID 12596)
$ws.Cells.Item(81, 8).Value = 7218.2
$ws.Cells.Item(81, 9).Value = 1147
$ws.Cells.Item(81, 10).Value = 16325
$ws.Cells.Item(81, 11).Value = 2294
$ws.Cells.Item(81, 12).Value = 32650
$ws.Cells.Item(81, 13).Value = -1233
$ws.Cells.Item(81, 14).Value = -34772
# Row 84 (Leve Item ID 12596)
$ws.Cells.Item(84, 8).Value = 7218.2
$ws.Cells.Item(84, 9).Value = 1147
$ws.Cells.Item(84, 10).Value = 16325
$ws.Cells.Item(84, 11).Value = 11470
$ws.Cells.Item(84, 12).Value = 163250
$ws.Cells.Item(84, 13).Value = -6166
$ws.Cells.Item(84, 14).Value = -173858
# Row 112 (Leve Item ID 25836)
$ws.Cells.Item(112, 8).Value = 43277.168
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 43277.168
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 43277.168
$ws.Cells.Item(112, 14).Value = -46231.168
